# v2.6 Added decoupled suspension, four-wheel steering, scripts to generate GGV diagram
#
# This script adds a new "FSAE_Achilles" aero-coefficients worksheet (cloned
# from the existing "Trailer_Kumanzi" sheet, which carries the shared
# layout/styles/tab-color used by every vehicle sheet in this workbook),
# fills in its class-specific values, tweaks one coefficient on the
# "Sedan_Hamba" sheet, and updates the active-sheet/selection bookkeeping to
# match where the author left off editing.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Add the new "FSAE_Achilles" vehicle sheet at the end of the tab strip,
#    cloned from "Trailer_Kumanzi" (sheetId 7) so it inherits the common
#    table layout, column widths, number formats and tab color.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("Trailer_Kumanzi")
$template.Copy([Type]::Missing, $template) | Out-Null
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "FSAE_Achilles"

# Fill in the FSAE_Achilles-specific values (class/instance name + aero
# coefficients), leaving the rest of the cloned template untouched.
$newSheet.Range("H3").Value = "FSAE_Achilles"
$newSheet.Range("H5").Value = -2.5
$newSheet.Range("H6").Value = 1
$newSheet.Range("H8").Value = 1.2
$newSheet.Range("F9").Value = -0.8
$newSheet.Range("G9").Value = 0
$newSheet.Range("H9").Value = 0.6

# Record the in-sheet selection the author left on this new tab.
$newSheet.Activate() | Out-Null
$newSheet.Range("G12").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) Sedan_Hamba: CD coefficient (H8) revised from 2.81 to 1.98, now shown
#    with a 2-decimal-place number format instead of General.
# ---------------------------------------------------------------------
$sedanHamba = $wb.Worksheets.Item("Sedan_Hamba")
$sedanHamba.Range("H8").Value = 1.98
$sedanHamba.Range("H8").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# 3) Selection bookkeeping: Sedan_HambaLG's last selected cell moved from
#    H6 to H8.
# ---------------------------------------------------------------------
$sedanHambaLG = $wb.Worksheets.Item("Sedan_HambaLG")
$sedanHambaLG.Activate() | Out-Null
$sedanHambaLG.Range("H8").Select() | Out-Null

# ---------------------------------------------------------------------
# 4) Sedan_Hamba becomes the active tab/selection (was Trailer_Kumanzi),
#    with its last selected cell now E18.
# ---------------------------------------------------------------------
$sedanHamba.Activate() | Out-Null
$sedanHamba.Range("E18").Select() | Out-Null
